$wb = $excel.ActiveWorkbook

# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# values for the ab719be0-... row (row 4) on both the zh-cn and de-de sheets,
# as part of generating the handback report.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-17 09:33:45"
$wsZhCn.Range("G4").Value = "2016-02-17 09:34:43"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-17 09:33:59"
$wsDeDe.Range("G4").Value = "2016-02-17 09:35:04"
